# Applies the betexplorer scrape-script update for greece/super-league-2 2023-2024
# Commit: "Atualizado por script em 05-11-2023 20:45"
#
# The underlying scraper re-ran and matches got reshuffled within a handful of
# row-groups (rotating the home/away/odds/url columns F:V among the rows of the
# group) plus one brand-new fixture appended as row 66 (Kozani FC vs AEL Larissa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("F6").Value = "AEL Larissa"
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = "Kampaniakos"
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1.29
$ws.Range("K6").Value = "23/09/2023 04:12"
$ws.Range("L6").Value = 1.31
$ws.Range("M6").Value = "23/09/2023 13:23"
$ws.Range("N6").Value = 4.51
$ws.Range("O6").Value = "23/09/2023 04:12"
$ws.Range("P6").Value = 4.77
$ws.Range("Q6").Value = "24/09/2023 13:03"
$ws.Range("R6").Value = 8.710000000000001
$ws.Range("S6").Value = "23/09/2023 04:12"
$ws.Range("T6").Value = 10.02
$ws.Range("U6").Value = "23/09/2023 13:23"
$ws.Range("V6").Value = "https://www.betexplorer.com/football/greece/super-league-2/ael-larissa-kampaniakos/2DC3RPYt/"

# Row 7
$ws.Range("F7").Value = "Tilikratis L."
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "PAE Egaleo"
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2.47
$ws.Range("K7").Value = "23/09/2023 04:12"
$ws.Range("L7").Value = 3.26
$ws.Range("M7").Value = "24/09/2023 14:31"
$ws.Range("N7").Value = 2.82
$ws.Range("O7").Value = "23/09/2023 04:12"
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = "24/09/2023 13:02"
$ws.Range("R7").Value = 2.82
$ws.Range("S7").Value = "23/09/2023 04:12"
$ws.Range("T7").Value = 2.32
$ws.Range("U7").Value = "24/09/2023 14:31"
$ws.Range("V7").Value = "https://www.betexplorer.com/football/greece/super-league-2/tilikratis-lefkada-pae-egaleo/pd8OM1mC/"

# Row 8
$ws.Range("F8").Value = "Panathinaikos B"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "Kalamata"
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 3.59
$ws.Range("K8").Value = "23/09/2023 03:13"
$ws.Range("L8").Value = 4.39
$ws.Range("M8").Value = "24/09/2023 14:57"
$ws.Range("N8").Value = 3.03
$ws.Range("O8").Value = "23/09/2023 03:13"
$ws.Range("P8").Value = 3.3
$ws.Range("Q8").Value = "24/09/2023 14:57"
$ws.Range("R8").Value = 1.97
$ws.Range("S8").Value = "23/09/2023 03:13"
$ws.Range("T8").Value = 1.85
$ws.Range("U8").Value = "24/09/2023 14:57"
$ws.Range("V8").Value = "https://www.betexplorer.com/football/greece/super-league-2/panathinaikos-kalamata/Wh4KNsY5/"

# Row 9
$ws.Range("F9").Value = "Karditsa"
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = "AEK Athens FC B"
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1.7
$ws.Range("K9").Value = "23/09/2023 03:13"
$ws.Range("L9").Value = 2.42
$ws.Range("M9").Value = "24/09/2023 14:42"
$ws.Range("N9").Value = 3.31
$ws.Range("O9").Value = "23/09/2023 03:13"
$ws.Range("P9").Value = 3.11
$ws.Range("Q9").Value = "24/09/2023 14:41"
$ws.Range("R9").Value = 4.43
$ws.Range("S9").Value = "23/09/2023 03:13"
$ws.Range("T9").Value = 2.97
$ws.Range("U9").Value = "24/09/2023 14:42"
$ws.Range("V9").Value = "https://www.betexplorer.com/football/greece/super-league-2/karditsa-aek/bsoWqYl8/"

# Row 34
$ws.Range("F34").Value = "PAOK B"
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = "Aiolikos"
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1.57
$ws.Range("K34").Value = "21/10/2023 01:12"
$ws.Range("L34").Value = 1.75
$ws.Range("M34").Value = "21/10/2023 14:43"
$ws.Range("N34").Value = 3.91
$ws.Range("O34").Value = "21/10/2023 01:12"
$ws.Range("P34").Value = 3.66
$ws.Range("Q34").Value = "21/10/2023 14:43"
$ws.Range("R34").Value = 5.06
$ws.Range("S34").Value = "21/10/2023 01:12"
$ws.Range("T34").Value = 4.44
$ws.Range("U34").Value = "21/10/2023 14:43"
$ws.Range("V34").Value = "https://www.betexplorer.com/football/greece/super-league-2/paok-aiolikos-fc/fyiX48y2/"

# Row 35
$ws.Range("F35").Value = "Athens Kallithea"
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = "Giouchtas"
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 1.3
$ws.Range("K35").Value = "21/10/2023 01:12"
$ws.Range("L35").Value = 1.23
$ws.Range("M35").Value = "21/10/2023 14:17"
$ws.Range("N35").Value = 4.55
$ws.Range("O35").Value = "21/10/2023 01:12"
$ws.Range("P35").Value = 5.44
$ws.Range("Q35").Value = "21/10/2023 14:18"
$ws.Range("R35").Value = 10.53
$ws.Range("S35").Value = "21/10/2023 01:12"
$ws.Range("T35").Value = 13.52
$ws.Range("U35").Value = "21/10/2023 14:17"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/greece/super-league-2/athens-kallithea-giouchtas/n5zBsZvU/"

# Row 36
$ws.Range("F36").Value = "Ionikos"
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = "Tilikratis L."
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 1.3
$ws.Range("K36").Value = "21/10/2023 01:12"
$ws.Range("L36").Value = 1.22
$ws.Range("M36").Value = "21/10/2023 12:52"
$ws.Range("N36").Value = 4.55
$ws.Range("O36").Value = "21/10/2023 01:12"
$ws.Range("P36").Value = 5.42
$ws.Range("Q36").Value = "21/10/2023 13:02"
$ws.Range("R36").Value = 10.53
$ws.Range("S36").Value = "21/10/2023 01:12"
$ws.Range("T36").Value = 14.23
$ws.Range("U36").Value = "21/10/2023 12:52"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/greece/super-league-2/ionikos-tilikratis-lefkada/hry7rFPN/"

# Row 39
$ws.Range("F39").Value = "Panachaiki"
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = "PAE Egaleo"
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 3.03
$ws.Range("K39").Value = "21/10/2023 02:13"
$ws.Range("L39").Value = 3.1
$ws.Range("M39").Value = "22/10/2023 14:59"
$ws.Range("N39").Value = 2.94
$ws.Range("O39").Value = "21/10/2023 02:13"
$ws.Range("P39").Value = 3.06
$ws.Range("Q39").Value = "22/10/2023 14:57"
$ws.Range("R39").Value = 2.3
$ws.Range("S39").Value = "21/10/2023 02:13"
$ws.Range("T39").Value = 2.37
$ws.Range("U39").Value = "22/10/2023 14:59"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/greece/super-league-2/panachaiki-pae-egaleo/QV5L1xQo/"

# Row 42
$ws.Range("F42").Value = "Kozani FC"
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = "Karditsa"
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2.75
$ws.Range("K42").Value = "22/10/2023 12:12"
$ws.Range("L42").Value = 2.58
$ws.Range("M42").Value = "22/10/2023 14:49"
$ws.Range("N42").Value = 2.74
$ws.Range("O42").Value = "22/10/2023 12:12"
$ws.Range("P42").Value = 2.74
$ws.Range("Q42").Value = "22/10/2023 13:54"
$ws.Range("R42").Value = 2.92
$ws.Range("S42").Value = "22/10/2023 12:12"
$ws.Range("T42").Value = 3.14
$ws.Range("U42").Value = "22/10/2023 14:49"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/greece/super-league-2/kozani-fc-karditsa/OnHhCjqS/"

# Row 46
$ws.Range("F46").Value = "Giouchtas"
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = "Panathinaikos B"
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1.83
$ws.Range("K46").Value = "28/10/2023 01:12"
$ws.Range("L46").Value = 1.91
$ws.Range("M46").Value = "28/10/2023 13:55"
$ws.Range("N46").Value = 3.35
$ws.Range("O46").Value = "28/10/2023 01:12"
$ws.Range("P46").Value = 3.24
$ws.Range("Q46").Value = "28/10/2023 13:55"
$ws.Range("R46").Value = 4.09
$ws.Range("S46").Value = "28/10/2023 01:12"
$ws.Range("T46").Value = 4.19
$ws.Range("U46").Value = "28/10/2023 13:55"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/greece/super-league-2/giouchtas-panathinaikos/AJ7XbzA4/"

# Row 50
$ws.Range("F50").Value = "Kalamata"
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = "Ilioupoli"
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 1.44
$ws.Range("K50").Value = "27/10/2023 02:13"
$ws.Range("L50").Value = 1.38
$ws.Range("M50").Value = "28/10/2023 13:43"
$ws.Range("N50").Value = 3.96
$ws.Range("O50").Value = "27/10/2023 02:13"
$ws.Range("P50").Value = 4.35
$ws.Range("Q50").Value = "28/10/2023 13:43"
$ws.Range("R50").Value = 5.97
$ws.Range("S50").Value = "27/10/2023 02:13"
$ws.Range("T50").Value = 8.56
$ws.Range("U50").Value = "28/10/2023 13:43"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/greece/super-league-2/kalamata-ilioupoli/bsP9un9q/"

# Row 63
$ws.Range("F63").Value = "Kampaniakos"
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = "Apollon Pontou"
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2.17
$ws.Range("K63").Value = "05/11/2023 03:13"
$ws.Range("L63").Value = 1.97
$ws.Range("M63").Value = "05/11/2023 12:16"
$ws.Range("N63").Value = 2.99
$ws.Range("O63").Value = "05/11/2023 03:13"
$ws.Range("P63").Value = 3.11
$ws.Range("Q63").Value = "05/11/2023 12:16"
$ws.Range("R63").Value = 3.42
$ws.Range("S63").Value = "05/11/2023 03:13"
$ws.Range("T63").Value = 4.12
$ws.Range("U63").Value = "05/11/2023 12:16"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/greece/super-league-2/kampaniakos-apollon-pontou/CUdC5uSC/"

# Row 64
$ws.Range("F64").Value = "Ionikos"
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = "Kalamata"
$ws.Range("I64").Value = 2
$ws.Range("J64").Value = 2.56
$ws.Range("K64").Value = "05/11/2023 03:13"
$ws.Range("L64").Value = 2.88
$ws.Range("M64").Value = "05/11/2023 13:59"
$ws.Range("N64").Value = 3
$ws.Range("O64").Value = "05/11/2023 03:13"
$ws.Range("P64").Value = 3.15
$ws.Range("Q64").Value = "05/11/2023 13:59"
$ws.Range("R64").Value = 2.82
$ws.Range("S64").Value = "05/11/2023 03:13"
$ws.Range("T64").Value = 2.46
$ws.Range("U64").Value = "05/11/2023 13:59"
$ws.Range("V64").Value = "https://www.betexplorer.com/football/greece/super-league-2/ionikos-kalamata/2JXbBSs2/"

# Row 65
$ws.Range("F65").Value = "Panathinaikos B"
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = "PAE Egaleo"
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2.48
$ws.Range("K65").Value = "04/11/2023 02:12"
$ws.Range("L65").Value = 2.59
$ws.Range("M65").Value = "05/11/2023 13:55"
$ws.Range("N65").Value = 2.88
$ws.Range("O65").Value = "04/11/2023 02:12"
$ws.Range("P65").Value = 3.01
$ws.Range("Q65").Value = "05/11/2023 13:55"
$ws.Range("R65").Value = 2.75
$ws.Range("S65").Value = "04/11/2023 02:12"
$ws.Range("T65").Value = 2.83
$ws.Range("U65").Value = "05/11/2023 03:41"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/greece/super-league-2/panathinaikos-pae-egaleo/OSYfC8Se/"

# Row 66 (new fixture appended)
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = "greece"
$ws.Range("C66").Value = "super-league-2"
$ws.Range("D66").Value = "2023-2024"
$ws.Range("E66").Value = 45235.58333333334
$ws.Range("F66").Value = "Kozani FC"
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = "AEL Larissa"
$ws.Range("I66").Value = 2
$ws.Range("J66").Value = 3.62
$ws.Range("K66").Value = "04/11/2023 02:12"
$ws.Range("L66").Value = 4.52
$ws.Range("M66").Value = "05/11/2023 13:40"
$ws.Range("N66").Value = 2.95
$ws.Range("O66").Value = "04/11/2023 02:12"
$ws.Range("P66").Value = 3.12
$ws.Range("Q66").Value = "05/11/2023 13:40"
$ws.Range("R66").Value = 2
$ws.Range("S66").Value = "04/11/2023 02:12"
$ws.Range("T66").Value = 1.88
$ws.Range("U66").Value = "05/11/2023 13:40"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/greece/super-league-2/kozani-fc-ael-larissa/U5776LC6/"

# Match the Indice/data_partida styling (bold+border+centered for A, date-time
# number format for E) used by every other data row, by copying formats from
# the first data row rather than re-declaring a brand-new style.
$ws.Range("A2").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E2").Copy()
$ws.Range("E66").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()
